# :sparkles: #2 컬럼 auto width 설정 옵션 추가
# Regenerate the random sample data (columns C, D, E) and the "now" timestamp
# columns (G, H) used by the column auto-width test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNow = 45020.69443405093

$data = @(
    @{ Row = 2;  C = 922319.0;  D = 415785.1875;  E = 975595.0 },
    @{ Row = 3;  C = 476448.0;  D = 652313.1875;  E = 742906.0 },
    @{ Row = 4;  C = 373542.0;  D = 430267.40625; E = 39273.0  },
    @{ Row = 5;  C = 864392.0;  D = 820477.75;    E = 896935.0 },
    @{ Row = 6;  C = 970491.0;  D = 725982.875;   E = 464732.0 },
    @{ Row = 7;  C = 734363.0;  D = 462756.9375;  E = 62625.0  },
    @{ Row = 8;  C = 823141.0;  D = 361631.03125; E = 245904.0 },
    @{ Row = 9;  C = 409575.0;  D = 562065.9375;  E = 629746.0 },
    @{ Row = 10; C = 914266.0;  D = 197526.8125;  E = 455287.0 },
    @{ Row = 11; C = 234904.0;  D = 928959.6875;  E = 546612.0 },
    @{ Row = 12; C = 792949.0;  D = 944007.25;    E = 865461.0 },
    @{ Row = 13; C = 509707.0;  D = 625906.25;    E = 917982.0 }
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 7).Value = $newNow
    $ws.Cells.Item($r, 8).Value = $newNow
}
